# Generate Report for Handoff
# Swap the old e2e markdown UUID / xliff hash for the newly generated ones,
# and bump the "Latest Handoff" / "Latest HO Xliff Generate" timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "4f0f6a26-27d9-4870-a4bc-d71dbbcd8276"
$newGuid = "0a24f56a-195f-4297-a987-da646670bcb0"

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c082a7a322c0a63d3455eb86e1cf2d9012520497/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name (A2), Path And Name (B2, hyperlinked),
# Latest HO Xliff Generate Date (G2)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$newB2Display = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkTarget, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newB2Display)

$wsOverview.Range("G2").Value = "2016-09-01 19:07:57"

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name (A2, hyperlinked), Latest Handoff File
# (G2), Latest Handoff Datetime (H2)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$newA2Display = "$newGuid.md"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkTarget, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newA2Display)

$wsZhCn.Range("G2").Value = "$newGuid.d72723ae0fc48b1b187dda5e4a8202c5d149647f.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-01 19:07:53"

# ---------------------------------------------------------------------
# Sheet "de-de": Source File Name (A2, hyperlinked), Latest Handoff File
# (G2), Latest Handoff Datetime (H2)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkTarget, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newA2Display)

$wsDeDe.Range("G2").Value = "$newGuid.d72723ae0fc48b1b187dda5e4a8202c5d149647f.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-01 19:07:57"
